{"js": "// Update the date line and the day's answer table with the new values.\nconst body = context.document.body;\n\n// 1) Update the date paragraph (first paragraph in the body).\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst datePara = paragraphs.items[0];\ndatePara.load(\"text\");\nawait context.sync();\n\nif (datePara.text === \"2025-11-27 Thursday\") {\n  datePara.insertText(\"2025-11-28 Friday\", \"Replace\");\n  await context.sync();\n}\n\n// 2) Update the answer table cell values (row-major 2D array),\n//    keeping the blank spacer rows untouched.\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\nconst oldToNew = new Map([\n  [\"40\u00f77=5, 5\", \"89\u00f72=44, 1\"],\n  [\"75\u00f77=10, 5\", \"60\u00f76=10, 0\"],\n  [\"82\u00f73=27, 1\", \"33\u00f79=3, 6\"],\n  [\"43\u00f72=21, 1\", \"44\u00f77=6, 2\"],\n  [\"69\u00f77=9, 6\", \"86\u00f72=43, 0\"],\n  [\"68\u00f74=17, 0\", \"67\u00f73=22, 1\"],\n  [\"89\u00f75=17, 4\", \"20\u00f74=5, 0\"],\n  [\"69\u00f75=13, 4\", \"98\u00f76=16, 2\"],\n  [\"21\u00f79=2, 3\", \"35\u00f77=5, 0\"],\n  [\"55\u00f72=27, 1\", \"91\u00f74=22, 3\"],\n  [\"60\u00f77=8, 4\", \"21\u00f77=3, 0\"],\n  [\"14\u00f77=2, 0\", \"89\u00f77=12, 5\"],\n  [\"19\u00f72=9, 1\", \"18\u00f73=6, 0\"],\n  [\"39\u00f76=6, 3\", \"13\u00f72=6, 1\"],\n  [\"34\u00f79=3, 7\", \"52\u00f76=8, 4\"],\n  [\"11\u00f78=1, 3\", \"33\u00f76=5, 3\"],\n  [\"64\u00f72=32, 0\", \"18\u00f74=4, 2\"],\n  [\"26\u00f73=8, 2\", \"41\u00f78=5, 1\"],\n  [\"34\u00f75=6, 4\", \"42\u00f74=10, 2\"],\n  [\"31\u00f78=3, 7\", \"22\u00f72=11, 0\"],\n  [\"51\u00f72=25, 1\", \"29\u00f78=3, 5\"],\n  [\"29\u00f72=14, 1\", \"73\u00f77=10, 3\"],\n  [\"43\u00f78=5, 3\", \"15\u00f78=1, 7\"],\n  [\"90\u00f77=12, 6\", \"54\u00f72=27, 0\"],\n  [\"56\u00f76=9, 2\", \"22\u00f72=11, 0\"],\n]);\n\nconst newValues = table.values.map((row) =>\n  row.map((cell) => (oldToNew.has(cell) ? oldToNew.get(cell) : cell))\n);\n\ntable.values = newValues;\nawait context.sync();\n", "ps1": "# Update the date line and the day's answer table with the new values.\n$d = $word.ActiveDocument\n\n# Ordered list of (find, replace) pairs -- the date paragraph followed by\n# every table cell value, in document order.\n$pairs = @(\n  @(\"2025-11-27 Thursday\", \"2025-11-28 Friday\"),\n  @(\"40\u00f77=5, 5\", \"89\u00f72=44, 1\"),\n  @(\"75\u00f77=10, 5\", \"60\u00f76=10, 0\"),\n  @(\"82\u00f73=27, 1\", \"33\u00f79=3, 6\"),\n  @(\"43\u00f72=21, 1\", \"44\u00f77=6, 2\"),\n  @(\"69\u00f77=9, 6\", \"86\u00f72=43, 0\"),\n  @(\"68\u00f74=17, 0\", \"67\u00f73=22, 1\"),\n  @(\"89\u00f75=17, 4\", \"20\u00f74=5, 0\"),\n  @(\"69\u00f75=13, 4\", \"98\u00f76=16, 2\"),\n  @(\"21\u00f79=2, 3\", \"35\u00f77=5, 0\"),\n  @(\"55\u00f72=27, 1\", \"91\u00f74=22, 3\"),\n  @(\"60\u00f77=8, 4\", \"21\u00f77=3, 0\"),\n  @(\"14\u00f77=2, 0\", \"89\u00f77=12, 5\"),\n  @(\"19\u00f72=9, 1\", \"18\u00f73=6, 0\"),\n  @(\"39\u00f76=6, 3\", \"13\u00f72=6, 1\"),\n  @(\"34\u00f79=3, 7\", \"52\u00f76=8, 4\"),\n  @(\"11\u00f78=1, 3\", \"33\u00f76=5, 3\"),\n  @(\"64\u00f72=32, 0\", \"18\u00f74=4, 2\"),\n  @(\"26\u00f73=8, 2\", \"41\u00f78=5, 1\"),\n  @(\"34\u00f75=6, 4\", \"42\u00f74=10, 2\"),\n  @(\"31\u00f78=3, 7\", \"22\u00f72=11, 0\"),\n  @(\"51\u00f72=25, 1\", \"29\u00f78=3, 5\"),\n  @(\"29\u00f72=14, 1\", \"73\u00f77=10, 3\"),\n  @(\"43\u00f78=5, 3\", \"15\u00f78=1, 7\"),\n  @(\"90\u00f77=12, 6\", \"54\u00f72=27, 0\"),\n  @(\"56\u00f76=9, 2\", \"22\u00f72=11, 0\")\n)\n\nforeach ($pair in $pairs) {\n  $findText = $pair[0]\n  $replaceText = $pair[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null\n}\n"}
